$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("A2").Value = "9KPTHTAQXE"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "44.5"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "08/04/2015"
$ws.Range("D2").Value = "IrregularIncome"
$ws.Range("E2").Value = "Eurofootball"

# Row 3
$ws.Range("A3").Value = "87VLWB1XEW"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "400"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "08/04/2015"
$ws.Range("D3").Value = "RegularIncome"
$ws.Range("E3").Value = "PayCheck"

# Row 4
$ws.Range("A4").Value = "5OVK26GTU5"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "10"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "08/04/2015"
$ws.Range("D4").Value = "IrregularExpense"
$ws.Range("E4").Value = "Lost"

# Restore the default (General, style 0) formatting on the text-coerced
# numeric/date-looking cells without disturbing their now-literal text
# content, by pasting just the format from an existing style-0 cell.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("B2:C4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
